$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.905.61'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').Value = '1.628.87'
$ws.Range('E3').Value = '  +1.59%  '
$ws.Range('E4').Value = '  +0.59%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('E8').Value = '  -4.31%  '
$ws.Range('E9').Value = '  +0.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0606'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0906'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').Value = '1.863.83'
$ws.Range('E12').Value = '  +1.69%  '
$ws.Range('D13').Value = '1.630.22'
$ws.Range('E13').Value = '  +1.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.560'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.98%  '
$ws.Range('D15').Value = '29.938.46'
$ws.Range('E15').Value = '  +0.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '9.07'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +13.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.83'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.94%  '
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '241.08'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('D20').Value = '0.0₃0700'
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('E21').Value = '  +0.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.11'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.70'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.82%  '
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.80'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.47'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.57'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.59%  '
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('E30').Value = '  +1.14%  '
$ws.Range('E31').Value = '  +3.38%  '
$ws.Range('E32').Value = '  +3.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.16'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.51%  '
$ws.Range('D34').Value = '1.422.71'
$ws.Range('E35').Value = '  +4.04%  '
$ws.Range('E36').Value = '  -1.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.76'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.87%  '
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '75.08'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +11.53%  '
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.98'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.824'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0490'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.82%  '
$ws.Range('E45').Value = '  +0.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '52.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.69%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '1.771.45'
$ws.Range('E48').Value = '  +1.79%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.33'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('E50').Value = '  +10.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '90.22'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.87%  '
